$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "maa://24702 (94.1), maa://25390 (96.6), maa://36681 (90.77)"
$ws.Range("AA2").Value = "maa://21246 (91.26), maa://36684 (98.63), ***maa://22731 (6.67)"
$ws.Range("K3").Value = "*maa://22880 (69.23), maa://20276 (82.86), *maa://22749 (66.67)"
$ws.Range("S4").Value = "maa://32509 (98.8), maa://22754 (91.67), maa://27295 (80.39), *maa://21746 (55.81), *maa://31008 (78.05)"
$ws.Range("C6").Value = "maa://42407 (88.89)"
$ws.Range("K6").Value = "maa://24839 (99.21)"
$ws.Range("S6").Value = "*maa://37411 (77.78)"
$ws.Range("O7").Value = "maa://22750 (97.14)"
$ws.Range("C8").Value = "*maa://21476 (69.77), **maa://39431 (50.0), *maa://37551 (57.14)"
$ws.Range("O8").Value = "maa://32931 (87.8), *maa://21916 (60.34), maa://23252 (92.31), **maa://22759 (45.45), maa://37496 (100.0)"
$ws.Range("W8").Value = "maa://21411 (96.07)"
$ws.Range("S9").Value = "**maa://22866 (30.77), maa://26222 (97.37)"
$ws.Range("S10").Value = "maa://27395 (96.08), maa://22755 (87.74), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("S11").Value = "maa://22747 (94.41), maa://22501 (98.15)"
$ws.Range("AA11").Value = "maa://22516 (89.29), maa://29912 (100.0), *maa://20794 (52.24)"
$ws.Range("G13").Value = "*maa://21248 (75.12), **maa://22728 (47.62)"
$ws.Range("O13").Value = "maa://22676 (91.84), *maa://22583 (75.41), *maa://22500 (55.81)"
$ws.Range("W13").Value = "*maa://34957 (76.6), *maa://22768 (51.61)"
$ws.Range("AE13").Value = "**maa://22737 (30.6), maa://39883 (88.89), *maa://39885 (73.68)"
$ws.Range("K14").Value = "maa://26245 (96.12), maa://21288 (96.21), maa://36682 (100.0), maa://39841 (93.33)"
$ws.Range("C15").Value = "*maa://22743 (76.88), maa://22734 (83.33), *maa://30808 (64.29), ***maa://36048 (12.9)"
$ws.Range("C16").Value = "maa://21441 (96.17), maa://36679 (91.43), maa://37650 (95.45)"
$ws.Range("O16").Value = "maa://28504 (91.84)"
$ws.Range("S16").Value = "maa://22729 (95.17), *maa://28648 (69.09), *maa://36674 (80.0)"
$ws.Range("AE16").Value = "*maa://23911 (62.37), maa://27755 (91.89)"
$ws.Range("C18").Value = "maa://24570 (96.67)"
$ws.Range("K21").Value = "maa://31731 (95.24)"
$ws.Range("AA21").Value = "*maa://21443 (78.79), **maa://23820 (30.91)"
$ws.Range("AE21").Value = "maa://22524 (94.32), *maa://22432 (74.55)"
$ws.Range("K22").Value = "maa://27127 (86.52), *maa://22751 (77.42)"
$ws.Range("W24").Value = "maa://29988 (86.41), maa://23504 (92.92), **maa://22892 (40.14), *maa://25141 (77.05), maa://36663 (80.7), ***maa://22815 (23.08)"
$ws.Range("AE24").Value = "maa://22523 (85.19), *maa://36672 (76.74), maa://29910 (94.12), **maa://21440 (34.55)"
$ws.Range("G25").Value = "*maa://29063 (75.56), *maa://25311 (74.19), ***maa://22725 (4.84)"
$ws.Range("G26").Value = "maa://24913 (91.18)"
$ws.Range("AA26").Value = "*maa://42235 (78.26)"
$ws.Range("C28").Value = "maa://24465 (90.37), maa://25725 (82.28)"
$ws.Range("W28").Value = "maa://39929 (86.89), ***maa://39723 (14.71), maa://41749 (82.35)"
$ws.Range("AE28").Value = "maa://36660 (93.89), *maa://36701 (64.0)"
$ws.Range("C29").Value = "maa://31694 (97.96)"
$ws.Range("K29").Value = "maa://28432 (93.54), *maa://28440 (72.84), maa://31400 (100.0), *maa://28650 (66.67)"
$ws.Range("G32").Value = "maa://21895 (97.01), maa://36667 (98.18), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("O33").Value = "*maa://21956 (79.23), maa://22730 (82.14)"
$ws.Range("S34").Value = "maa://24526 (93.19)"
$ws.Range("K35").Value = "maa://41296 (98.18)"
$ws.Range("O37").Value = "maa://21280 (89.19), *maa://21239 (72.73)"
$ws.Range("AE38").Value = "maa://36697 (84.21)"
$ws.Range("G39").Value = "maa://25199 (86.11), maa://36670 (88.24), maa://30434 (87.5), ***maa://25036 (16.0)"
$ws.Range("O40").Value = "maa://23278 (95.89), maa://21386 (95.65), maa://36664 (90.48)"
$ws.Range("O41").Value = "**maa://35616 (36.67)"
$ws.Range("G44").Value = "maa://29768 (97.55), maa://27728 (96.0)"
$ws.Range("G46").Value = "maa://35931 (92.64)"
$ws.Range("G53").Value = "maa://32534 (93.31), **maa://32434 (36.36)"
